{"js": "// Insert three new bullet paragraphs into the Siege Analytics (PARTNER)\n// section, right after the \"GIS & Geospatial Analysis Consulting\" line and\n// before the \"Lead comprehensive research...\" bullet.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Locate the anchor paragraph by its exact text.\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"GIS & Geospatial Analysis Consulting\") {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find anchor paragraph 'GIS & Geospatial Analysis Consulting'\");\n}\n\nconst newBullets = [\n  \"\u2022 Utilized ESRI Arc Suite and OSGeo technology to map and analyze 50,000+ electoral boundaries across federal, state, and local levels\",\n  \"\u2022 Applied geospatial analysis to uncover demographic miscoding affecting 2,000+ precincts nationwide\",\n  \"\u2022 Developed boundary estimation tools enabling smaller organizations to conduct sophisticated redistricting analysis\"\n];\n\n// Insert the three paragraphs directly after the anchor, preserving order\n// by always inserting immediately after the anchor itself (each new\n// paragraph becomes the new \"next\" sibling, so we chain off the anchor).\nlet insertAfter = anchor;\nfor (const text of newBullets) {\n  insertAfter = insertAfter.insertParagraph(text, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# Insert three new bullet paragraphs into the Siege Analytics (PARTNER)\n# section, right after the \"GIS & Geospatial Analysis Consulting\" line and\n# before the \"Lead comprehensive research...\" bullet.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph by its exact text.\n$anchorIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.TrimEnd(\"`r`a\") -eq \"GIS & Geospatial Analysis Consulting\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find anchor paragraph 'GIS & Geospatial Analysis Consulting'\"\n}\n\n$newBullets = @(\n    \"\u2022 Utilized ESRI Arc Suite and OSGeo technology to map and analyze 50,000+ electoral boundaries across federal, state, and local levels\",\n    \"\u2022 Applied geospatial analysis to uncover demographic miscoding affecting 2,000+ precincts nationwide\",\n    \"\u2022 Developed boundary estimation tools enabling smaller organizations to conduct sophisticated redistricting analysis\"\n)\n\n# Insert each bullet directly after the anchor paragraph, in order, so that\n# each subsequent insertion lands right after the previously inserted one.\n$currentIndex = $anchorIndex\nforeach ($text in $newBullets) {\n    $currentPara = $d.Paragraphs.Item($currentIndex)\n    $currentPara.Range.InsertParagraphAfter()\n    $currentIndex = $currentIndex + 1\n    $newPara = $d.Paragraphs.Item($currentIndex)\n    $newPara.Range.Text = $text\n}\n"}
